$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B3").Value = 1500.75
$ws.Range("B4").Value = 0.75
$ws.Range("B5").Value = 0.68
$ws.Range("B6").Value = 22
$ws.Range("B7").Value = 12
$ws.Range("B9").Value = 54.55

# --- Strategy Status sheet ---
$ws = $wb.Worksheets.Item("Strategy Status")
$ws.Range("C6").Value = 100.75
$ws.Range("D6").Value = 22
$ws.Range("E6").Value = 0.75
$ws.Range("F6").Value = 0.75
$ws.Range("G6").Value = 54.55

# --- New trade (row 23) appended to "All Trades" and "MarketMaking" sheets ---
function Set-TradeRow23($ws) {
    $ws.Cells.Item(23, 1).Value = 22

    $dateCell = $ws.Cells.Item(23, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"
    $dateCell.ClearFormats()

    $timeCell = $ws.Cells.Item(23, 3)
    $timeCell.NumberFormat = "@"
    $timeCell.Value = "23:57:08"
    $timeCell.ClearFormats()

    $ws.Cells.Item(23, 4).Value = "MarketMaking"
    $ws.Cells.Item(23, 5).Value = "UP"
    $ws.Cells.Item(23, 6).Value = 0.5600000000000001
    $ws.Cells.Item(23, 7).Value = 0.58
    $ws.Cells.Item(23, 8).Value = "CLOSED"
    $ws.Cells.Item(23, 9).Value = 3.5714
    $ws.Cells.Item(23, 10).Value = 0.02
    $ws.Cells.Item(23, 11).Value = 100.75
    $ws.Cells.Item(23, 12).Value = 0
    $ws.Cells.Item(23, 13).Value = 0
    $ws.Cells.Item(23, 14).Value = 0.6
    $ws.Cells.Item(23, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(23, 16).Value = "early_exit"
    $ws.Cells.Item(23, 17).Value = 0.14
}

Set-TradeRow23 $wb.Worksheets.Item("All Trades")
Set-TradeRow23 $wb.Worksheets.Item("MarketMaking")
